# Updated cryptos list values (Price / Volume(1h) columns) to match the
# refreshed coinranking.com snapshot. Cells whose new text parses as a
# plain number are briefly forced to Text format ("@") before the
# assignment (so e.g. "52.00" / "0.3633" keep their exact digits instead
# of being coerced to the numbers 52 / 0.3633...) and then restored to the
# sheet's default "Normal" style so no extra formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '23.222.79'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").Value = '1.602.37'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.94'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.86%  '
$ws.Range("E7").Value = '  -0.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '52.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3633'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.274'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08147'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.21%  '
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.79'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.578'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.407'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.02%  '
$ws.Range("E16").Value = '  +0.63%  '
$ws.Range("D17").Value = '1.599.36'
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06921'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.530'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.34%  '
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.91'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.33%  '
$ws.Range("D24").Value = '23.215.90'
$ws.Range("E24").Value = '  +0.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.454'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.041'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.22'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.87%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '149.92'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.281'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '135.87'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.385'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.737'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.48%  '
$ws.Range("D33").Value = '1.778.61'
$ws.Range("E33").Value = '  -0.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9647'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07491'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.35'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02736'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.24%  '
$ws.Range("E38").Value = '  -0.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.131'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08789'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.386'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7105'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.61'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6543'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.317'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.31%  '
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.008'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '132.55'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07926'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.205'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.33%  '
